$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update libraries_tools for the credit_supervised row (H8) with the longest
# string, replacing the old "pandas, sklearn, imblearn" value.
$ws.Range("H8").Value = "pandas, splinter, bs4 (beautiful soup 4), sqlalchemy, bootstrap, AWS, postgresql, flask "

# Update selection to match new active cell (H8)
$ws.Range("H8").Select()
